$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
# D1 used to be "status"; it becomes "obtained_result" and a new E1 "status" is added.
$ws.Range("D1").Value = "obtained_result"
$ws.Range("E1").Value = "status"

# --- Clear old D2:D11 "passed" values (no longer present in the data) ---
$ws.Range("D2:D11").ClearContents()

# --- Replace the data rows (A:C) with the new dataset ---
$data = @(
    @(882, 1, "Low-Volume Loan Short-Term"),
    @(928, 2, "Low-Volume Loan Short-Term"),
    @(138, 3, "Low-Volume Loan Short-Term"),
    @(824, 3, "Low-Volume Loan Short-Term"),
    @(35,  2, "Low-Volume Loan Short-Term"),
    @(244, 7, "Low-Volume Loan Long-Term"),
    @(290, 8, "Low-Volume Loan Long-Term"),
    @(140, 7, "Low-Volume Loan Long-Term"),
    @(827, 6, "Low-Volume Loan Long-Term"),
    @(678, 5, "Low-Volume Loan Long-Term")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# --- Column widths ---
# (Excel's ColumnWidth setter snaps to whole-pixel boundaries, so the inputs
# below are chosen to land on the closest achievable stored width.)
$ws.Columns.Item(4).ColumnWidth = 28
$ws.Columns.Item(5).ColumnWidth = 7.75

# --- Selection ---
$ws.Range("D2:E11").Select()
